# Generate Report for Handoff
# Inserts a new row (for the file d5ae9e09-953a-4365-84f5-54a4117f4756.md) above
# the existing f1f67407-... row on all three sheets (Overview, zh-cn, de-de),
# pushing the existing row down, then fixes up hyperlinks, table ranges and
# the first column width.

$wb = $excel.ActiveWorkbook

$newBase = "d5ae9e09-953a-4365-84f5-54a4117f4756"
$oldBase = "f1f67407-6312-4f4f-8163-988cffa4f170"

$newMd = $newBase + ".md"
$oldMd = $oldBase + ".md"

$newUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/772338ad293faee52853fd9b8488311cf555d197/e2e/" + $newMd
$oldUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/772338ad293faee52853fd9b8488311cf555d197/e2e/" + $oldMd

$newDisplay = "e2e\" + $newMd
$oldDisplay = "e2e\" + $oldMd

# Width value that round-trips to an xlsx <col> width of exactly 40.
$colAWidth = 39.166666666666664

# ---------------------------------------------------------------------
# Sheet "Overview" (File Name / Path And Name / Extension / Publish URL /
# zh-cn / de-de / Latest HO Xliff Generate Date). Hyperlink lives in col B.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows.Item(3).Insert()

# Row 3 = old row, shifted down (copy old row2 values down first).
$wsOverview.Range("A3").Value = $oldMd
$wsOverview.Range("B3").Value = $oldDisplay
$wsOverview.Range("C3").Value = $wsOverview.Range("C2").Value()
$wsOverview.Range("E3").Value = $wsOverview.Range("E2").Value()
$wsOverview.Range("F3").Value = $wsOverview.Range("F2").Value()
$wsOverview.Range("G3").Value = "2016-08-24 02:39:30"

# Row 2 = new file's data.
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newDisplay
$wsOverview.Range("G2").Value = "2016-08-24 02:39:47"

# Hyperlinks: col B, rows 2 and 3.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B3").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $newUrl, "", "", $newDisplay)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $oldUrl, "", "", $oldDisplay)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Columns.Item(1).ColumnWidth = $colAWidth

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de" share the same 16-column layout:
# A Source File Name, B File Extension, C Status, D Source Path,
# E Priority, F Content Duplicate, G Latest Handoff File,
# H Latest Handoff Datetime, I Latest Target File, J Latest Handback File,
# K Latest Handback DateTime, L Reference Tokens, M To be localized,
# N Dependency From, O Has metadata, P Error Detail.
# Hyperlink lives in col A.
# ---------------------------------------------------------------------
function Update-LangSheet($ws, $newXlfSuffix, $newDateTime, $oldXlfSuffix, $oldDateTime) {
    $ws.Rows.Item(3).Insert()

    # Row 3 = old row, shifted down (copy old row2 values that are unaffected).
    $ws.Range("A3").Value = $oldMd
    $ws.Range("B3").Value = $ws.Range("B2").Value()
    $ws.Range("C3").Value = $ws.Range("C2").Value()
    $ws.Range("D3").Value = $ws.Range("D2").Value()
    $ws.Range("E3").Value = $ws.Range("E2").Value()
    $ws.Range("F3").Value = $ws.Range("F2").Value()
    $ws.Range("G3").Value = $oldBase + "." + $oldXlfSuffix
    $ws.Range("H3").Value = $oldDateTime
    $ws.Range("K3").Value = $ws.Range("K2").Value()
    $ws.Range("M3").Value = $ws.Range("M2").Value()
    $ws.Range("O3").Value = $ws.Range("O2").Value()

    # Row 2 = new file's data.
    $ws.Range("A2").Value = $newMd
    $ws.Range("G2").Value = $newBase + "." + $newXlfSuffix
    $ws.Range("H2").Value = $newDateTime

    # Hyperlinks: col A, rows 2 and 3.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("A3").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $newUrl, "", "", $newMd)
    $ws.Hyperlinks.Add($ws.Range("A3"), $oldUrl, "", "", $oldMd)

    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A1:P3"))

    $ws.Columns.Item(1).ColumnWidth = $colAWidth
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZhCn "44ef045c4a424fd02424c14b48a05e638c8cae54.zh-cn.xlf" "2016-08-24 02:39:42" "28dd6556fab734355a0e38485da294cb08e132fd.zh-cn.xlf" "2016-08-24 02:39:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDeDe "44ef045c4a424fd02424c14b48a05e638c8cae54.de-de.xlf" "2016-08-24 02:39:47" "28dd6556fab734355a0e38485da294cb08e132fd.de-de.xlf" "2016-08-24 02:39:30"
